# Update the "取得日時" (acquisition timestamp) column for the data rows
# on the active sheet ("ランサーズ") from 2025-09-26 18:23:25 to
# 2025-09-26 18:29:51, reflecting a new append run at 18:29 JST.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-26 18:29:51"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
